$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table S6 part subsets scores")

$ws.Range("A2").Value = "long COVID"
$ws.Range("B2").Value = "Tyrol"
$ws.Range("C2").Value = "Sum of acute symptoms"
$ws.Range("D2").Value = "mean(SD) = 17 (5.95)`nmedian(IQR) = 17 (13 - 21)`nrange = 4 - 34`nncomplete = 270"
$ws.Range("E2").Value = "mean(SD) = 15.8 (6.88)`nmedian(IQR) = 14.5 (11 - 20)`nrange = 2 - 41`nncomplete = 104"
$ws.Range("F2").Value = "mean(SD) = 16.4 (6.24)`nmedian(IQR) = 16 (12 - 20)`nrange = 2 - 33`nncomplete = 173"
$ws.Range("G2").Value = "ns (p = 0.11)"
$ws.Range("H2").Value = "ns (p = 0.28)"

$ws.Range("A3").Value = "long COVID"
$ws.Range("B3").Value = "Tyrol"
$ws.Range("C3").Value = "Physical performance loss"
$ws.Range("D3").Value = "mean(SD) = 26.9 (21.5)`nmedian(IQR) = 20 (10 - 39)`nrange = 0 - 100`nncomplete = 269"
$ws.Range("E3").Value = "mean(SD) = 21 (19.7)`nmedian(IQR) = 15 (7.75 - 30)`nrange = 0 - 89`nncomplete = 104"
$ws.Range("F3").Value = "mean(SD) = 22.6 (20.2)`nmedian(IQR) = 19 (10 - 31)`nrange = 0 - 92`nncomplete = 173"
$ws.Range("G3").Value = "p = 0.01"
$ws.Range("H3").Value = "ns (p = 0.05)"

$ws.Range("A4").Value = "long COVID"
$ws.Range("B4").Value = "Tyrol"
$ws.Range("C4").Value = "Quality of life score"
$ws.Range("D4").Value = "mean(SD) = 1.17 (0.78)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 270"
$ws.Range("E4").Value = "mean(SD) = 1.16 (0.739)`nmedian(IQR) = 1 (1 - 1.25)`nrange = 0 - 3`nncomplete = 104"
$ws.Range("F4").Value = "mean(SD) = 1.1 (0.805)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 173"
$ws.Range("G4").Value = "ns (p = 0.69)"
$ws.Range("H4").Value = "ns (p = 0.69)"

$ws.Range("A5").Value = "long COVID"
$ws.Range("B5").Value = "Tyrol"
$ws.Range("C5").Value = "Overall mental health score"
$ws.Range("D5").Value = "mean(SD) = 1.2 (0.777)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 270"
$ws.Range("E5").Value = "mean(SD) = 1.08 (0.784)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 104"
$ws.Range("F5").Value = "mean(SD) = 1.08 (0.75)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 173"
$ws.Range("G5").Value = "ns (p = 0.2)"
$ws.Range("H5").Value = "ns (p = 0.34)"

$ws.Range("A6").Value = "long COVID"
$ws.Range("B6").Value = "Tyrol"
$ws.Range("C6").Value = "PHQ stress score"
$ws.Range("D6").Value = "mean(SD) = 5.37 (3.59)`nmedian(IQR) = 5 (3 - 7)`nrange = 0 - 18`nncomplete = 269"
$ws.Range("E6").Value = "mean(SD) = 5.24 (3.8)`nmedian(IQR) = 4 (2.75 - 8)`nrange = 0 - 16`nncomplete = 104"
$ws.Range("F6").Value = "mean(SD) = 5.03 (3.9)`nmedian(IQR) = 4 (2 - 7)`nrange = 0 - 19`nncomplete = 172"
$ws.Range("G6").Value = "ns (p = 0.33)"
$ws.Range("H6").Value = "ns (p = 0.41)"

$ws.Range("A7").Value = "long COVID"
$ws.Range("B7").Value = "South Tyrol"
$ws.Range("C7").Value = "Sum of acute symptoms"
$ws.Range("D7").Value = "mean(SD) = 17.6 (6.68)`nmedian(IQR) = 17 (13 - 22)`nrange = 4 - 38`nncomplete = 244"
$ws.Range("E7").Value = "mean(SD) = 13.5 (6.76)`nmedian(IQR) = 14 (7.5 - 17)`nrange = 2 - 34`nncomplete = 59"
$ws.Range("F7").Value = "mean(SD) = 17.7 (7.89)`nmedian(IQR) = 16.5 (12 - 24)`nrange = 3 - 39`nncomplete = 136"
$ws.Range("G7").Value = "p = 0.00052"
$ws.Range("H7").Value = "p = 0.0013"

$ws.Range("A8").Value = "long COVID"
$ws.Range("B8").Value = "South Tyrol"
$ws.Range("C8").Value = "Physical performance loss"
$ws.Range("D8").Value = "mean(SD) = 25.4 (20.2)`nmedian(IQR) = 20 (10 - 34)`nrange = 0 - 93`nncomplete = 240"
$ws.Range("E8").Value = "mean(SD) = 14.8 (16.8)`nmedian(IQR) = 10 (2.5 - 20)`nrange = 0 - 75`nncomplete = 58"
$ws.Range("F8").Value = "mean(SD) = 26.2 (23.8)`nmedian(IQR) = 20 (8 - 40)`nrange = 0 - 100`nncomplete = 136"
$ws.Range("G8").Value = "p = 0.00033"
$ws.Range("H8").Value = "p = 0.0013"

$ws.Range("A9").Value = "long COVID"
$ws.Range("B9").Value = "South Tyrol"
$ws.Range("C9").Value = "Quality of life score"
$ws.Range("D9").Value = "mean(SD) = 1.29 (0.697)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 244"
$ws.Range("E9").Value = "mean(SD) = 1.07 (0.666)`nmedian(IQR) = 1 (1 - 1)`nrange = 0 - 3`nncomplete = 59"
$ws.Range("F9").Value = "mean(SD) = 1.29 (0.807)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 136"
$ws.Range("G9").Value = "ns (p = 0.096)"
$ws.Range("H9").Value = "ns (p = 0.096)"

$ws.Range("A10").Value = "long COVID"
$ws.Range("B10").Value = "South Tyrol"
$ws.Range("C10").Value = "Overall mental health score"
$ws.Range("D10").Value = "mean(SD) = 1.22 (0.779)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 244"
$ws.Range("E10").Value = "mean(SD) = 0.915 (0.702)`nmedian(IQR) = 1 (0 - 1)`nrange = 0 - 3`nncomplete = 59"
$ws.Range("F10").Value = "mean(SD) = 1.2 (0.806)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 136"
$ws.Range("G10").Value = "p = 0.022"
$ws.Range("H10").Value = "p = 0.027"

$ws.Range("A11").Value = "long COVID"
$ws.Range("B11").Value = "South Tyrol"
$ws.Range("C11").Value = "PHQ stress score"
$ws.Range("D11").Value = "mean(SD) = 5.11 (3.29)`nmedian(IQR) = 5 (2 - 7)`nrange = 0 - 15`nncomplete = 242"
$ws.Range("E11").Value = "mean(SD) = 3.61 (2.7)`nmedian(IQR) = 3 (1.5 - 6)`nrange = 0 - 13`nncomplete = 59"
$ws.Range("F11").Value = "mean(SD) = 5.68 (4.11)`nmedian(IQR) = 5 (2.75 - 8)`nrange = 0 - 19`nncomplete = 136"
$ws.Range("G11").Value = "p = 0.0021"
$ws.Range("H11").Value = "p = 0.0035"

$ws.Range("A12").Value = "PASC"
$ws.Range("B12").Value = "Tyrol"
$ws.Range("C12").Value = "Sum of acute symptoms"
$ws.Range("D12").Value = "mean(SD) = 17.9 (7.17)`nmedian(IQR) = 18 (13 - 22)`nrange = 5 - 42`nncomplete = 103"
$ws.Range("E12").Value = "mean(SD) = 15.6 (6.65)`nmedian(IQR) = 14 (11 - 20)`nrange = 2 - 33`nncomplete = 54"
$ws.Range("F12").Value = "mean(SD) = 16.2 (5.96)`nmedian(IQR) = 17 (13 - 20)`nrange = 5 - 29`nncomplete = 57"
$ws.Range("G12").Value = "ns (p = 0.11)"
$ws.Range("H12").Value = "ns (p = 0.14)"

$ws.Range("A13").Value = "PASC"
$ws.Range("B13").Value = "Tyrol"
$ws.Range("C13").Value = "Physical performance loss"
$ws.Range("D13").Value = "mean(SD) = 29.7 (22)`nmedian(IQR) = 24 (11.2 - 45.8)`nrange = 0 - 100`nncomplete = 102"
$ws.Range("E13").Value = "mean(SD) = 20.3 (20)`nmedian(IQR) = 14 (8.25 - 30)`nrange = 0 - 92`nncomplete = 54"
$ws.Range("F13").Value = "mean(SD) = 20.9 (18.7)`nmedian(IQR) = 19 (10 - 30)`nrange = 0 - 92`nncomplete = 57"
$ws.Range("G13").Value = "p = 0.003"
$ws.Range("H13").Value = "p = 0.015"

$ws.Range("A14").Value = "PASC"
$ws.Range("B14").Value = "Tyrol"
$ws.Range("C14").Value = "Quality of life score"
$ws.Range("D14").Value = "mean(SD) = 1.26 (0.766)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 103"
$ws.Range("E14").Value = "mean(SD) = 1.07 (0.64)`nmedian(IQR) = 1 (1 - 1)`nrange = 0 - 3`nncomplete = 54"
$ws.Range("F14").Value = "mean(SD) = 1.07 (0.821)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 57"
$ws.Range("G14").Value = "ns (p = 0.19)"
$ws.Range("H14").Value = "ns (p = 0.19)"

$ws.Range("A15").Value = "PASC"
$ws.Range("B15").Value = "Tyrol"
$ws.Range("C15").Value = "Overall mental health score"
$ws.Range("D15").Value = "mean(SD) = 1.3 (0.826)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 103"
$ws.Range("E15").Value = "mean(SD) = 1.04 (0.8)`nmedian(IQR) = 1 (0.25 - 1.75)`nrange = 0 - 3`nncomplete = 54"
$ws.Range("F15").Value = "mean(SD) = 1.05 (0.692)`nmedian(IQR) = 1 (1 - 1)`nrange = 0 - 3`nncomplete = 57"
$ws.Range("G15").Value = "ns (p = 0.065)"
$ws.Range("H15").Value = "ns (p = 0.11)"

$ws.Range("A16").Value = "PASC"
$ws.Range("B16").Value = "Tyrol"
$ws.Range("C16").Value = "PHQ stress score"
$ws.Range("D16").Value = "mean(SD) = 6.25 (4.15)`nmedian(IQR) = 5 (3 - 9.5)`nrange = 0 - 16`nncomplete = 103"
$ws.Range("E16").Value = "mean(SD) = 4.94 (3.8)`nmedian(IQR) = 4 (3 - 6.75)`nrange = 0 - 16`nncomplete = 54"
$ws.Range("F16").Value = "mean(SD) = 4.72 (3.89)`nmedian(IQR) = 3 (2 - 7)`nrange = 0 - 19`nncomplete = 57"
$ws.Range("G16").Value = "p = 0.027"
$ws.Range("H16").Value = "ns (p = 0.067)"

$ws.Range("A17").Value = "PASC"
$ws.Range("B17").Value = "South Tyrol"
$ws.Range("C17").Value = "Sum of acute symptoms"
$ws.Range("D17").Value = "mean(SD) = 18.6 (7.37)`nmedian(IQR) = 18 (14 - 23)`nrange = 2 - 38`nncomplete = 98"
$ws.Range("E17").Value = "mean(SD) = 13.3 (6.64)`nmedian(IQR) = 15 (7 - 17)`nrange = 5 - 25`nncomplete = 21"
$ws.Range("F17").Value = "mean(SD) = 17.8 (8.21)`nmedian(IQR) = 16 (11.8 - 25)`nrange = 4 - 33`nncomplete = 52"
$ws.Range("G17").Value = "p = 0.022"
$ws.Range("H17").Value = "p = 0.043"

$ws.Range("A18").Value = "PASC"
$ws.Range("B18").Value = "South Tyrol"
$ws.Range("C18").Value = "Physical performance loss"
$ws.Range("D18").Value = "mean(SD) = 28.7 (22.4)`nmedian(IQR) = 24 (11 - 45.8)`nrange = 0 - 93`nncomplete = 98"
$ws.Range("E18").Value = "mean(SD) = 15.4 (16.5)`nmedian(IQR) = 12 (3.25 - 23)`nrange = 0 - 58`nncomplete = 20"
$ws.Range("F18").Value = "mean(SD) = 24.4 (20.2)`nmedian(IQR) = 20 (10 - 38.2)`nrange = 0 - 75`nncomplete = 52"
$ws.Range("G18").Value = "p = 0.026"
$ws.Range("H18").Value = "p = 0.043"

$ws.Range("A19").Value = "PASC"
$ws.Range("B19").Value = "South Tyrol"
$ws.Range("C19").Value = "Quality of life score"
$ws.Range("D19").Value = "mean(SD) = 1.43 (0.718)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 98"
$ws.Range("E19").Value = "mean(SD) = 0.905 (0.7)`nmedian(IQR) = 1 (0 - 1)`nrange = 0 - 2`nncomplete = 21"
$ws.Range("F19").Value = "mean(SD) = 1.33 (0.734)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 52"
$ws.Range("G19").Value = "p = 0.017"
$ws.Range("H19").Value = "p = 0.043"

$ws.Range("A20").Value = "PASC"
$ws.Range("B20").Value = "South Tyrol"
$ws.Range("C20").Value = "Overall mental health score"
$ws.Range("D20").Value = "mean(SD) = 1.24 (0.787)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 98"
$ws.Range("E20").Value = "mean(SD) = 1.05 (0.74)`nmedian(IQR) = 1 (1 - 1)`nrange = 0 - 3`nncomplete = 21"
$ws.Range("F20").Value = "mean(SD) = 1.19 (0.817)`nmedian(IQR) = 1 (1 - 2)`nrange = 0 - 3`nncomplete = 52"
$ws.Range("G20").Value = "ns (p = 0.52)"
$ws.Range("H20").Value = "ns (p = 0.52)"

$ws.Range("A21").Value = "PASC"
$ws.Range("B21").Value = "South Tyrol"
$ws.Range("C21").Value = "PHQ stress score"
$ws.Range("D21").Value = "mean(SD) = 5.77 (3.29)`nmedian(IQR) = 6 (3 - 8)`nrange = 0 - 15`nncomplete = 98"
$ws.Range("E21").Value = "mean(SD) = 3.9 (1.92)`nmedian(IQR) = 4 (2 - 6)`nrange = 1 - 7`nncomplete = 21"
$ws.Range("F21").Value = "mean(SD) = 5.96 (3.91)`nmedian(IQR) = 5.5 (2.75 - 8.25)`nrange = 0 - 14`nncomplete = 52"
$ws.Range("G21").Value = "p = 0.047"
$ws.Range("H21").Value = "ns (p = 0.059)"
